$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 533.5
$ws.Range("I5").Value = 65.5
$ws.Range("J5").Value = 1001.5
$ws.Range("K5").Value = 65.5
$ws.Range("L5").Value = 1001.5
$ws.Range("M5").Value = 49.5
$ws.Range("N5").Value = -1231.5
$ws.Range("H12").Value = 145.33333
$ws.Range("I12").Value = 145.33333
$ws.Range("K12").Value = 145.33333
$ws.Range("M12").Value = 24.66667000000001
$ws.Range("H86").Value = 191319.4
$ws.Range("I86").Value = 1201
$ws.Range("J86").Value = 318065
$ws.Range("K86").Value = 1201
$ws.Range("L86").Value = 318065
$ws.Range("M86").Value = -78
$ws.Range("N86").Value = -320311
$ws.Range("H89").Value = 191319.4
$ws.Range("I89").Value = 1201
$ws.Range("J89").Value = 318065
$ws.Range("K89").Value = 6005
$ws.Range("L89").Value = 1590325
$ws.Range("M89").Value = -389
$ws.Range("N89").Value = -1601557
$ws.Range("H113").Value = 3533
$ws.Range("I113").Value = 1599
$ws.Range("K113").Value = 1599
$ws.Range("M113").Value = 1655
$ws.Range("H129").Value = 3788.5
$ws.Range("I129").Value = 4797.3335
$ws.Range("K129").Value = 14392.0005
$ws.Range("M129").Value = -9392.000499999998
$ws.Range("H137").Value = 1299.1666
$ws.Range("I137").Value = 1198.75
$ws.Range("K137").Value = 3596.25
$ws.Range("M137").Value = -1046.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1991.762
$ws.Range("I2").Value = 1557.909
$ws.Range("K2").Value = 1557.909
$ws.Range("M2").Value = -1444.909
$ws.Range("H32").Value = 3211489.8
$ws.Range("J32").Value = 2333998.2
$ws.Range("L32").Value = 2333998.2
$ws.Range("N32").Value = -2334572.2
$ws.Range("H45").Value = 4933.1665
$ws.Range("I45").Value = 5697.8
$ws.Range("J45").Value = 1110
$ws.Range("K45").Value = 5697.8
$ws.Range("L45").Value = 1110
$ws.Range("M45").Value = -5320.8
$ws.Range("N45").Value = -1864
$ws.Range("H50").Value = 4797.5713
$ws.Range("I50").Value = 9473.333000000001
$ws.Range("K50").Value = 9473.333000000001
$ws.Range("M50").Value = -8759.333000000001
$ws.Range("H110").Value = 1671
$ws.Range("I110").Value = 1808.1
$ws.Range("K110").Value = 1808.1
$ws.Range("M110").Value = 236.9000000000001
$ws.Range("H116").Value = 1991.762
$ws.Range("I116").Value = 1557.909
$ws.Range("K116").Value = 1557.909
$ws.Range("M116").Value = 736.0909999999999
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").Value = ""

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1991.762
$ws.Range("I3").Value = 1557.909
$ws.Range("K3").Value = 1557.909
$ws.Range("M3").Value = -1443.909
$ws.Range("H86").Value = 1193.3684
$ws.Range("I86").Value = 1033.9166
$ws.Range("K86").Value = 1033.9166
$ws.Range("M86").Value = 89.08339999999998
$ws.Range("H89").Value = 1193.3684
$ws.Range("I89").Value = 1033.9166
$ws.Range("K89").Value = 5169.583000000001
$ws.Range("M89").Value = 446.4169999999995
$ws.Range("H94").Value = 1596.15
$ws.Range("I94").Value = 1440.2222
$ws.Range("K94").Value = 1440.2222
$ws.Range("M94").Value = -989.2221999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 813.5
$ws.Range("I22").Value = 767
$ws.Range("J22").Value = 860
$ws.Range("K22").Value = 767
$ws.Range("L22").Value = 860
$ws.Range("M22").Value = -417
$ws.Range("N22").Value = -1560
$ws.Range("H31").Value = 2098.5
$ws.Range("I31").Value = 2098
$ws.Range("K31").Value = 2098
$ws.Range("M31").Value = -1803
$ws.Range("H34").Value = 2098.5
$ws.Range("I34").Value = 2098
$ws.Range("K34").Value = 2098
$ws.Range("M34").Value = -1896
$ws.Range("H99").Value = 1824.6666
$ws.Range("J99").Value = 2090.3333
$ws.Range("L99").Value = 2090.3333
$ws.Range("N99").Value = -5086.3333
$ws.Range("H107").Value = 848.1111
$ws.Range("I107").Value = 733.2857
$ws.Range("K107").Value = 733.2857
$ws.Range("M107").Value = 1186.7143
$ws.Range("H126").Value = 1824.6666
$ws.Range("J126").Value = 2090.3333
$ws.Range("L126").Value = 6270.999899999999
$ws.Range("N126").Value = -11210.9999
$ws.Range("H134").Value = 2495.4707
$ws.Range("I134").Value = 2428.2
$ws.Range("K134").Value = 7284.599999999999
$ws.Range("M134").Value = -4749.599999999999
$ws.Range("H141").Value = 180871.14
$ws.Range("J141").Value = 180871.14
$ws.Range("L141").Value = 180871.14
$ws.Range("N141").Value = -191231.14

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 183.5
$ws.Range("I12").Value = 271.83334
$ws.Range("K12").Value = 815.5000200000001
$ws.Range("M12").Value = -642.5000200000001
$ws.Range("H132").Value = 51500
$ws.Range("J132").Value = 3000
$ws.Range("L132").Value = 27000
$ws.Range("N132").Value = -32060
$ws.Range("H137").Value = 3455.4
$ws.Range("I137").Value = 2392.5
$ws.Range("J137").Value = 4164
$ws.Range("K137").Value = 7177.5
$ws.Range("L137").Value = 12492
$ws.Range("M137").Value = -2077.5
$ws.Range("N137").Value = -22692

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 57.142857
$ws.Range("I2").Value = 48.2
$ws.Range("J2").Value = 79.5
$ws.Range("K2").Value = 48.2
$ws.Range("L2").Value = 79.5
$ws.Range("M2").Value = 64.8
$ws.Range("N2").Value = -305.5
$ws.Range("H11").Value = 2175000.5
$ws.Range("J11").Value = 757250
$ws.Range("L11").Value = 757250
$ws.Range("N11").Value = -757528
$ws.Range("H70").Value = 5008
$ws.Range("I70").Value = 5008
$ws.Range("K70").Value = 5008
$ws.Range("M70").Value = -4738
$ws.Range("H73").Value = 5008
$ws.Range("I73").Value = 5008
$ws.Range("K73").Value = 5008
$ws.Range("M73").Value = -4072
$ws.Range("H107").Value = 1293.6666
$ws.Range("J107").Value = 2948
$ws.Range("L107").Value = 2948
$ws.Range("N107").Value = -6788
$ws.Range("H122").Value = 4095.8
$ws.Range("I122").Value = 4095.8
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 12287.4
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -9837.400000000001
$ws.Range("N122").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7551.4443
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").Value = ""
$ws.Range("H22").Value = 1706.8572
$ws.Range("I22").Value = 1824.75
$ws.Range("J22").Value = 1549.6666
$ws.Range("K22").Value = 1824.75
$ws.Range("L22").Value = 1549.6666
$ws.Range("M22").Value = -1529.75
$ws.Range("N22").Value = -2139.6666
$ws.Range("H27").Value = 1706.8572
$ws.Range("I27").Value = 1824.75
$ws.Range("J27").Value = 1549.6666
$ws.Range("K27").Value = 1824.75
$ws.Range("L27").Value = 1549.6666
$ws.Range("M27").Value = -1717.75
$ws.Range("N27").Value = -1763.6666
$ws.Range("H40").Value = 3752.3125
$ws.Range("I40").Value = 3476.4546
$ws.Range("K40").Value = 3476.4546
$ws.Range("M40").Value = -3340.4546
$ws.Range("H46").Value = 2060.647
$ws.Range("I46").Value = 2179
$ws.Range("J46").Value = 1996.091
$ws.Range("K46").Value = 2179
$ws.Range("L46").Value = 1996.091
$ws.Range("M46").Value = -1991
$ws.Range("N46").Value = -2372.091
$ws.Range("H55").Value = 1401.9375
$ws.Range("I55").Value = 1366.8572
$ws.Range("J55").Value = 1429.2222
$ws.Range("K55").Value = 1366.8572
$ws.Range("L55").Value = 1429.2222
$ws.Range("M55").Value = -1193.8572
$ws.Range("N55").Value = -1775.2222
$ws.Range("H126").Value = 7551.4443
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 40002
$ws.Range("I62").Value = 40002
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 40002
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -39378
$ws.Range("N62").Value = ""
$ws.Range("H65").Value = 40002
$ws.Range("I65").Value = 40002
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 200010
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -196890
$ws.Range("N65").Value = ""
$ws.Range("H95").Value = 21625.334
$ws.Range("J95").Value = 21625.334
$ws.Range("L95").Value = 21625.334
$ws.Range("N95").Value = -27117.334
$ws.Range("H122").Value = 9318.333000000001
$ws.Range("I122").Value = 9857.333000000001
$ws.Range("J122").Value = 8779.333000000001
$ws.Range("K122").Value = 29571.999
$ws.Range("L122").Value = 26337.999
$ws.Range("M122").Value = -27121.999
$ws.Range("N122").Value = -31237.999
